# DSS.xlsx update: add two new trainees (Abdallah Kasem Awad Mahmoud,
# Ismail Abdulaal Ismail Attia Elmelegy) with their 8 course rows each,
# restyle the now-contiguous "FAWZY ABDELKADER ALI ZAHRA" block, and move
# the active selection down to the newly entered data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) B462:B469 ("FAWZY ABDELKADER ALI ZAHRA") loses its special fill and
#    now matches the plain style used by the rest of that row block.
# ---------------------------------------------------------------------
$ws.Range("A462").Copy()
$ws.Range("B462:B469").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Rows 470-477: Abdallah Kasem Awad Mahmoud's eight courses.
#    Copy formatting from the row block directly above (454:461), which
#    already carries the exact style combination these rows need.
# ---------------------------------------------------------------------
$ws.Range("A454:E461").Copy()
$ws.Range("A470:E477").PasteSpecial(-4122)

$ws.Range("A470").Value = "DSS1469"
$ws.Range("B470").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C470").Value = "30 Hours Construction Safety & Health"
$ws.Range("D470").Value = "'04-11-2024"
$ws.Range("E470").Value = 1

$ws.Range("A471").Value = "DSS1470"
$ws.Range("B471").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C471").Value = "30 Hours G. Industry Safety & Health"
$ws.Range("D471").Value = "'09-11-2024"
$ws.Range("E471").Value = 1

$ws.Range("A472").Value = "DSS1471"
$ws.Range("B472").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C472").Value = "Electrical Safety "
$ws.Range("D472").Value = "'06-11-2024"
$ws.Range("E472").Value = 1

$ws.Range("A473").Value = "DSS1472"
$ws.Range("B473").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C473").Value = "Fire Marshal"
$ws.Range("D473").Value = "'03-11-2024"
$ws.Range("E473").Value = 1

$ws.Range("A474").Value = "DSS1473"
$ws.Range("B474").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C474").Value = "Scaffold Competent Person"
$ws.Range("D474").Value = "'01-11-2024"
$ws.Range("E474").Value = 1

$ws.Range("A475").Value = "DSS1474"
$ws.Range("B475").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C475").Value = "Lifting & Rigging Competent Person"
$ws.Range("D475").Value = "'02-11-2024"
$ws.Range("E475").Value = 1

$ws.Range("A476").Value = "DSS1475"
$ws.Range("B476").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C476").Value = "Health & Safety Risk Assessment"
$ws.Range("D476").Value = "'07-11-2024"
$ws.Range("E476").Value = 1

$ws.Range("A477").Value = "DSS1476"
$ws.Range("B477").Value = "Abdallah Kasem Awad Mahmoud"
$ws.Range("C477").Value = "Safety Management System & PTW"
$ws.Range("D477").Value = "'08-11-2024"
$ws.Range("E477").Value = 1

# ---------------------------------------------------------------------
# 3) Rows 478-485: Ismail Abdulaal Ismail Attia Elmelegy's eight courses.
#    Copy formatting from the freshly-restyled 462:469 block (now
#    s=10/10/10/38/20), which is exactly the style this block needs.
# ---------------------------------------------------------------------
$ws.Range("A462:E469").Copy()
$ws.Range("A478:E485").PasteSpecial(-4122)

$ws.Range("A478").Value = "DSS1477"
$ws.Range("B478").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C478").Value = "30 Hours Construction Safety & Health"
$ws.Range("D478").Value = "'04-11-2024"
$ws.Range("E478").Value = 1

$ws.Range("A479").Value = "DSS1478"
$ws.Range("B479").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C479").Value = "30 Hours G. Industry Safety & Health"
$ws.Range("D479").Value = "'09-11-2024"
$ws.Range("E479").Value = 1

$ws.Range("A480").Value = "DSS1479"
$ws.Range("B480").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C480").Value = "Electrical Safety "
$ws.Range("D480").Value = "'06-11-2024"
$ws.Range("E480").Value = 1

$ws.Range("A481").Value = "DSS1480"
$ws.Range("B481").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C481").Value = "Fire Marshal"
$ws.Range("D481").Value = "'03-11-2024"
$ws.Range("E481").Value = 1

$ws.Range("A482").Value = "DSS1481"
$ws.Range("B482").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C482").Value = "Scaffold Competent Person"
$ws.Range("D482").Value = "'01-11-2024"
$ws.Range("E482").Value = 1

$ws.Range("A483").Value = "DSS1482"
$ws.Range("B483").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C483").Value = "Lifting & Rigging Competent Person"
$ws.Range("D483").Value = "'02-11-2024"
$ws.Range("E483").Value = 1

$ws.Range("A484").Value = "DSS1483"
$ws.Range("B484").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C484").Value = "Health & Safety Risk Assessment"
$ws.Range("D484").Value = "'07-11-2024"
$ws.Range("E484").Value = 1

$ws.Range("A485").Value = "DSS1484"
$ws.Range("B485").Value = "Ismail Abdulaal Ismail Attia Elmelegy"
$ws.Range("C485").Value = "Safety Management System & PTW"
$ws.Range("D485").Value = "'08-11-2024"
$ws.Range("E485").Value = 1

# ---------------------------------------------------------------------
# 4) Move the view down to where the new rows were entered.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 471
$ws.Range("C490").Select()
